# Auto-generated edit script
# Applies numeric value corrections to the Leve profit-tracking tables
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR),
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ----- Worksheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1434.875
$ws.Range("I2").Value = 1434.875
$ws.Range("K2").Value = 1434.875
$ws.Range("M2").Value = -1321.875
$ws.Range("H19").Value = 3126
$ws.Range("I19").Value = 195.33333
$ws.Range("K19").Value = 195.33333
$ws.Range("M19").Value = -20.33332999999999
$ws.Range("H21").Value = 33006.332
$ws.Range("I21").Value = 33006.332
$ws.Range("K21").Value = 33006.332
$ws.Range("M21").Value = -32538.332
$ws.Range("H23").Value = 33006.332
$ws.Range("I23").Value = 33006.332
$ws.Range("K23").Value = 33006.332
$ws.Range("M23").Value = -32772.332
$ws.Range("H28").Value = 595
$ws.Range("I28").Value = 133.06667
$ws.Range("J28").Value = 1749.8334
$ws.Range("K28").Value = 133.06667
$ws.Range("L28").Value = 1749.8334
$ws.Range("M28").Value = 351.93333
$ws.Range("N28").Value = -2719.8334
$ws.Range("H29").Value = 3000.6667
$ws.Range("I29").Value = 3000.6667
$ws.Range("K29").Value = 9002.000100000001
$ws.Range("M29").Value = -8721.000100000001
$ws.Range("H38").Value = 1240.4166
$ws.Range("I38").Value = 190.84616
$ws.Range("J38").Value = 2480.818
$ws.Range("K38").Value = 572.5384799999999
$ws.Range("L38").Value = 7442.454000000001
$ws.Range("M38").Value = -200.5384799999999
$ws.Range("N38").Value = -8186.454000000001
$ws.Range("H39").Value = 712.6667
$ws.Range("J39").Value = 1417.8
$ws.Range("L39").Value = 4253.4
$ws.Range("N39").Value = -4845.4
$ws.Range("H43").Value = 1740.037
$ws.Range("I43").Value = 1000.25
$ws.Range("J43").Value = 1868.6957
$ws.Range("K43").Value = 1000.25
$ws.Range("L43").Value = 1868.6957
$ws.Range("M43").Value = -931.25
$ws.Range("N43").Value = -2006.6957
$ws.Range("H58").Value = 1308.1818
$ws.Range("I58").Value = 198.18182
$ws.Range("J58").Value = 2418.182
$ws.Range("K58").Value = 594.5454599999999
$ws.Range("L58").Value = 7254.545999999999
$ws.Range("M58").Value = -444.5454599999999
$ws.Range("N58").Value = -7554.545999999999
$ws.Range("H62").Value = 3470.8125
$ws.Range("J62").Value = 4507.5
$ws.Range("L62").Value = 4507.5
$ws.Range("N62").Value = -5755.5
$ws.Range("H65").Value = 3470.8125
$ws.Range("J65").Value = 4507.5
$ws.Range("L65").Value = 22537.5
$ws.Range("N65").Value = -28777.5
$ws.Range("H107").Value = 1703.95
$ws.Range("I107").Value = 1048.5
$ws.Range("J107").Value = 3233.3333
$ws.Range("K107").Value = 1048.5
$ws.Range("L107").Value = 3233.3333
$ws.Range("M107").Value = 871.5
$ws.Range("N107").Value = -7073.3333
$ws.Range("H116").Value = 3034.2632
$ws.Range("I116").Value = 2545.9048
$ws.Range("J116").Value = 3637.5293
$ws.Range("K116").Value = 2545.9048
$ws.Range("L116").Value = 3637.5293
$ws.Range("M116").Value = 896.0952000000002
$ws.Range("N116").Value = -10521.5293
$ws.Range("H137").Value = 2328023.5
$ws.Range("I137").Value = 3573948.2
$ws.Range("J137").Value = 2297.3333
$ws.Range("K137").Value = 10721844.6
$ws.Range("L137").Value = 6891.999899999999
$ws.Range("M137").Value = -10719294.6
$ws.Range("N137").Value = -11991.9999

# ----- Worksheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4639.704
$ws.Range("I32").Value = 3550.4138
$ws.Range("K32").Value = 3550.4138
$ws.Range("M32").Value = -3263.4138
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314

# ----- Worksheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1760.84
$ws.Range("I99").Value = 1295
$ws.Range("K99").Value = 1295
$ws.Range("M99").Value = 203

# ----- Worksheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 16766.2
$ws.Range("I10").Value = 529.2
$ws.Range("J10").Value = 33003.2
$ws.Range("K10").Value = 529.2
$ws.Range("L10").Value = 33003.2
$ws.Range("M10").Value = -390.2
$ws.Range("N10").Value = -33281.2
$ws.Range("H31").Value = 1300767.6
$ws.Range("I31").Value = 2440339.8
$ws.Range("J31").Value = 2921.5278
$ws.Range("K31").Value = 2440339.8
$ws.Range("L31").Value = 2921.5278
$ws.Range("M31").Value = -2440044.8
$ws.Range("N31").Value = -3511.5278
$ws.Range("H34").Value = 1300767.6
$ws.Range("I34").Value = 2440339.8
$ws.Range("J34").Value = 2921.5278
$ws.Range("K34").Value = 2440339.8
$ws.Range("L34").Value = 2921.5278
$ws.Range("M34").Value = -2440137.8
$ws.Range("N34").Value = -3325.5278
$ws.Range("H99").Value = 1989.3334
$ws.Range("I99").Value = 1258.2727
$ws.Range("J99").Value = 3999.75
$ws.Range("K99").Value = 1258.2727
$ws.Range("L99").Value = 3999.75
$ws.Range("M99").Value = 239.7273
$ws.Range("N99").Value = -6995.75
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H126").Value = 1989.3334
$ws.Range("I126").Value = 1258.2727
$ws.Range("J126").Value = 3999.75
$ws.Range("K126").Value = 3774.8181
$ws.Range("L126").Value = 11999.25
$ws.Range("M126").Value = -1304.8181
$ws.Range("N126").Value = -16939.25
$ws.Range("H134").Value = 1745.711
$ws.Range("I134").Value = 1408.3024
$ws.Range("K134").Value = 4224.9072
$ws.Range("M134").Value = -1689.9072

# ----- Worksheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 812.1613
$ws.Range("I5").Value = 469.7037
$ws.Range("K5").Value = 1409.1111
$ws.Range("M5").Value = -1297.1111
$ws.Range("H46").Value = 2184
$ws.Range("I46").Value = 1125
$ws.Range("J46").Value = 2448.75
$ws.Range("K46").Value = 3375
$ws.Range("L46").Value = 7346.25
$ws.Range("M46").Value = -3284
$ws.Range("N46").Value = -7528.25
$ws.Range("H68").Value = 2147.224
$ws.Range("I68").Value = 678.5862
$ws.Range("K68").Value = 2035.7586
$ws.Range("M68").Value = -1224.7586
$ws.Range("H71").Value = 2147.224
$ws.Range("I71").Value = 678.5862
$ws.Range("K71").Value = 6107.275799999999
$ws.Range("M71").Value = -2051.275799999999
$ws.Range("H94").Value = 3682.4138
$ws.Range("I94").Value = 1633.3334
$ws.Range("J94").Value = 3918.8462
$ws.Range("K94").Value = 4900.0002
$ws.Range("L94").Value = 11756.5386
$ws.Range("M94").Value = -4224.0002
$ws.Range("N94").Value = -13108.5386
$ws.Range("H132").Value = 2923.0952
$ws.Range("I132").Value = 2488.2222
$ws.Range("J132").Value = 3249.25
$ws.Range("K132").Value = 22393.9998
$ws.Range("L132").Value = 29243.25
$ws.Range("M132").Value = -19863.9998
$ws.Range("N132").Value = -34303.25
$ws.Range("H135").Value = 812.1613
$ws.Range("I135").Value = 469.7037
$ws.Range("K135").Value = 4227.3333
$ws.Range("M135").Value = -1692.3333

# ----- Worksheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2875.1667
$ws.Range("I7").Value = 1060.8
$ws.Range("J7").Value = 4171.143
$ws.Range("K7").Value = 1060.8
$ws.Range("L7").Value = 4171.143
$ws.Range("M7").Value = -948.8
$ws.Range("N7").Value = -4395.143
$ws.Range("H40").Value = 3142.8572
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3142.8572
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3142.8572
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3414.8572
$ws.Range("H45").Value = 6985.25
$ws.Range("J45").Value = 8400
$ws.Range("L45").Value = 8400
$ws.Range("N45").Value = -9214
$ws.Range("H100").Value = 2808.5715
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 3110
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 3110
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -4192
$ws.Range("H126").Value = 2875.1667
$ws.Range("I126").Value = 1060.8
$ws.Range("J126").Value = 4171.143
$ws.Range("K126").Value = 3182.4
$ws.Range("L126").Value = 12513.429
$ws.Range("M126").Value = -712.3999999999996
$ws.Range("N126").Value = -17453.429
$ws.Range("H132").Value = 2780.5
$ws.Range("I132").Value = 1673.8823
$ws.Range("J132").Value = 4490.727
$ws.Range("K132").Value = 5021.6469
$ws.Range("L132").Value = 13472.181
$ws.Range("M132").Value = -2491.6469
$ws.Range("N132").Value = -18532.181
$ws.Range("H136").Value = 2705360.8
$ws.Range("I136").Value = 3705938.8
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 11117816.4
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -11115266.4
$ws.Range("N136").Value = -16500

# ----- Worksheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 324664.06
$ws.Range("I122").Value = 455913
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 1367739
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").Value = -1365289
$ws.Range("N122").Value = -16399.9999
$ws.Range("H136").Value = 1575.7949
$ws.Range("I136").Value = 692.40625
$ws.Range("J136").Value = 5614.143
$ws.Range("K136").Value = 2077.21875
$ws.Range("L136").Value = 16842.429
$ws.Range("M136").Value = 472.78125
$ws.Range("N136").Value = -21942.429

Write-Host "Applied all Leve profit updates."
